$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INVERTER")
$ws.Columns("L:M").Insert()

$ws.Range("L1").Value = "ACkW"
$ws.Range("M1").Value = "ACkVar"

$ws.Range("L2").Value = "'50"
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value = "'100"
$ws.Range("M2").Style = "Normal"

$ws.Range("L3").Value = "'0"
$ws.Range("L3").Style = "Normal"
$ws.Range("M3").Value = "'0"
$ws.Range("M3").Style = "Normal"

$ws.Range("L4").Value = "'0"
$ws.Range("L4").Style = "Normal"
$ws.Range("M4").Value = "'0"
$ws.Range("M4").Style = "Normal"

Write-Host "Done"
